$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A21").Value = "Centraal Station_B"
$ws.Range("D6").Select()
